$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows where Target cluster = "MuSCs" (old rows 16, 12, 8, 4), deleting
# from the bottom up so row numbers of earlier rows stay valid.
$ws.Rows("16").Delete()
$ws.Rows("12").Delete()
$ws.Rows("8").Delete()
$ws.Rows("4").Delete()

# Refresh the TPM-derived metrics for the remaining 12 sender/target combinations.
# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il15"
$ws.Range("C2").Value = "Il2rb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.314527000000001
$ws.Range("H2").Value = 18.943581
$ws.Range("I2").Value = 0.2616724966426195
$ws.Range("J2").Value = 0.2616724966426195
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1082756666666667
$ws.Range("N2").Value = 0.324827
$ws.Range("O2").Value = 0.05996835676729485
$ws.Range("P2").Value = 0.05996835676729485
$ws.Range("Q2").Value = 0.6837096206096667
$ws.Range("R2").Value = 6.153386585487
$ws.Range("S2").Value = 0.01569206963485337
$ws.Range("T2").Value = 0.01569206963485337

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il15"
$ws.Range("C3").Value = "Il2rb"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.314527000000001
$ws.Range("H3").Value = 18.943581
$ws.Range("I3").Value = 0.2616724966426195
$ws.Range("J3").Value = 0.2616724966426195
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.89265
$ws.Range("N3").Value = 2.67795
$ws.Range("O3").Value = 0.4943932031665386
$ws.Range("P3").Value = 0.4943932031665387
$ws.Range("Q3").Value = 5.636662526550001
$ws.Range("R3").Value = 50.72996273895001
$ws.Range("S3").Value = 0.12936910379573
$ws.Range("T3").Value = 0.12936910379573

# Row 4: ECs -> Resolving-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il15"
$ws.Range("C4").Value = "Il2rb"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.314527000000001
$ws.Range("H4").Value = 18.943581
$ws.Range("I4").Value = 0.2616724966426195
$ws.Range("J4").Value = 0.2616724966426195
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.804621
$ws.Range("N4").Value = 2.413863
$ws.Range("O4").Value = 0.4456384400661664
$ws.Range("P4").Value = 0.4456384400661665
$ws.Range("Q4").Value = 5.080801029267001
$ws.Range("R4").Value = 45.727209263403
$ws.Range("S4").Value = 0.1166113232120361
$ws.Range("T4").Value = 0.1166113232120362

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il15"
$ws.Range("C5").Value = "Il2rb"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.525638333333333
$ws.Range("H5").Value = 10.576915
$ws.Range("I5").Value = 0.1461016137776048
$ws.Range("J5").Value = 0.1461016137776048
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1082756666666667
$ws.Range("N5").Value = 0.324827
$ws.Range("O5").Value = 0.05996835676729485
$ws.Range("P5").Value = 0.05996835676729485
$ws.Range("Q5").Value = 0.3817408409672222
$ws.Range("R5").Value = 3.435667568705
$ws.Range("S5").Value = 0.008761473699292923
$ws.Range("T5").Value = 0.008761473699292923

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il15"
$ws.Range("C6").Value = "Il2rb"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.525638333333333
$ws.Range("H6").Value = 10.576915
$ws.Range("I6").Value = 0.1461016137776048
$ws.Range("J6").Value = 0.1461016137776048
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.89265
$ws.Range("N6").Value = 2.67795
$ws.Range("O6").Value = 0.4943932031665386
$ws.Range("P6").Value = 0.4943932031665387
$ws.Range("Q6").Value = 3.14716105825
$ws.Range("R6").Value = 28.32444952425
$ws.Range("S6").Value = 0.07223164482331051
$ws.Range("T6").Value = 0.07223164482331053

# Row 7: FAPs -> Resolving-Mac
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il15"
$ws.Range("C7").Value = "Il2rb"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.525638333333333
$ws.Range("H7").Value = 10.576915
$ws.Range("I7").Value = 0.1461016137776048
$ws.Range("J7").Value = 0.1461016137776048
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.804621
$ws.Range("N7").Value = 2.413863
$ws.Range("O7").Value = 0.4456384400661664
$ws.Range("P7").Value = 0.4456384400661665
$ws.Range("Q7").Value = 2.836802641405
$ws.Range("R7").Value = 25.531223772645
$ws.Range("S7").Value = 0.06510849525500133
$ws.Range("T7").Value = 0.06510849525500133

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Il15"
$ws.Range("C8").Value = "Il2rb"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.510814
$ws.Range("H8").Value = 1.532442
$ws.Range("I8").Value = 0.02116801063642661
$ws.Range("J8").Value = 0.02116801063642662
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1082756666666667
$ws.Range("N8").Value = 0.324827
$ws.Range("O8").Value = 0.05996835676729485
$ws.Range("P8").Value = 0.05996835676729485
$ws.Range("Q8").Value = 0.05530872639266666
$ws.Range("R8").Value = 0.497778537534
$ws.Range("S8").Value = 0.001269410813899123
$ws.Range("T8").Value = 0.001269410813899123

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Il15"
$ws.Range("C9").Value = "Il2rb"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.510814
$ws.Range("H9").Value = 1.532442
$ws.Range("I9").Value = 0.02116801063642661
$ws.Range("J9").Value = 0.02116801063642662
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.89265
$ws.Range("N9").Value = 2.67795
$ws.Range("O9").Value = 0.4943932031665386
$ws.Range("P9").Value = 0.4943932031665387
$ws.Range("Q9").Value = 0.4559781171
$ws.Range("R9").Value = 4.1038030539
$ws.Range("S9").Value = 0.01046532058320631
$ws.Range("T9").Value = 0.01046532058320632

# Row 10: MuSCs -> Resolving-Mac
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Il15"
$ws.Range("C10").Value = "Il2rb"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.510814
$ws.Range("H10").Value = 1.532442
$ws.Range("I10").Value = 0.02116801063642661
$ws.Range("J10").Value = 0.02116801063642662
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.804621
$ws.Range("N10").Value = 2.413863
$ws.Range("O10").Value = 0.4456384400661664
$ws.Range("P10").Value = 0.4456384400661665
$ws.Range("Q10").Value = 0.411011671494
$ws.Range("R10").Value = 3.699105043446
$ws.Range("S10").Value = 0.009433279239321175
$ws.Range("T10").Value = 0.009433279239321178

# Row 11: Resolving-Mac -> ECs
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Il15"
$ws.Range("C11").Value = "Il2rb"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 13.78043333333333
$ws.Range("H11").Value = 41.3413
$ws.Range("I11").Value = 0.571057878943349
$ws.Range("J11").Value = 0.5710578789433491
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1082756666666667
$ws.Range("N11").Value = 0.324827
$ws.Range("O11").Value = 0.05996835676729485
$ws.Range("P11").Value = 0.05996835676729485
$ws.Range("Q11").Value = 1.492085606122222
$ws.Range("R11").Value = 13.4287704551
$ws.Range("S11").Value = 0.03424540261924942
$ws.Range("T11").Value = 0.03424540261924943

# Row 12: Resolving-Mac -> FAPs
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Il15"
$ws.Range("C12").Value = "Il2rb"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 13.78043333333333
$ws.Range("H12").Value = 41.3413
$ws.Range("I12").Value = 0.571057878943349
$ws.Range("J12").Value = 0.5710578789433491
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.89265
$ws.Range("N12").Value = 2.67795
$ws.Range("O12").Value = 0.4943932031665386
$ws.Range("P12").Value = 0.4943932031665387
$ws.Range("Q12").Value = 12.301103815
$ws.Range("R12").Value = 110.709934335
$ws.Range("S12").Value = 0.2823271339642918
$ws.Range("T12").Value = 0.2823271339642919

# Row 13: Resolving-Mac -> Resolving-Mac
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Il15"
$ws.Range("C13").Value = "Il2rb"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 13.78043333333333
$ws.Range("H13").Value = 41.3413
$ws.Range("I13").Value = 0.571057878943349
$ws.Range("J13").Value = 0.5710578789433491
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.804621
$ws.Range("N13").Value = 2.413863
$ws.Range("O13").Value = 0.4456384400661664
$ws.Range("P13").Value = 0.4456384400661665
$ws.Range("Q13").Value = 11.0880260491
$ws.Range("R13").Value = 99.7922344419
$ws.Range("S13").Value = 0.2544853423598078
$ws.Range("T13").Value = 0.2544853423598079
